$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "42.994.07"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +2.42%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.307.44"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +2.08%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.11"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.84%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.55"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.536"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.11"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.16%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0818"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("E12").Value = "  +0.78%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.11"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +7.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.665.36"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.95"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +4.50%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.312.43"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +1.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.803"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.94%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "42.967.60"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +2.59%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  +2.50%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.08"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.22%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "68.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +0.97%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "240.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.32%  "
$ws.Range("E24").Value = "  +5.05%  "
$ws.Range("E25").Value = "  +2.22%  "
$ws.Range("E26").Value = "  +0.31%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "24.13"
$ws.Range("D27").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.96"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +6.37%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.67"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.08%  "
$ws.Range("E30").Value = "  +1.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "168.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +6.12%  "
$ws.Range("E32").Value = "  +2.88%  "
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.14"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "17.72"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.76%  "
$ws.Range("E36").Value = "  +0.66%  "
$ws.Range("E37").Value = "  +0.25%  "
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("E40").Value = "  +1.98%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.25"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +6.22%  "
$ws.Range("E42").Value = "  -5.01%  "
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "19.38"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +2.86%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.972.53"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0289"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.26%  "
$ws.Range("E46").Value = "  +3.49%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  +18.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "55.20"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.98%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.532.98"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.88%  "
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +2.94%  "
